$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dSF (column F) values per repulled data / mean calculation
$ws.Range("F2").Value = -2
$ws.Range("F3").Value = -2
$ws.Range("F4").Value = -4
$ws.Range("F6").Value = -1
$ws.Range("F7").Value = 2
$ws.Range("F8").Value = -2
